# freeCrmTestData.xlsx — "added pages & tests"
#   * Sheet2 renamed to "Companies" and populated with company test data
#   * Companies sheet becomes the active tab / selected sheet
#   * Contacts sheet selection moves from C4 to header row A1:D1

$wb = $excel.ActiveWorkbook

# --- Contacts sheet: drop the old single-cell selection, select the header row instead
$contacts = $wb.Worksheets.Item("Contacts")
$contacts.Range("A1:D1").Select() | Out-Null

# --- Sheet2 -> Companies
$companies = $wb.Worksheets.Item("Sheet2")
$companies.Name = "Companies"

$headers = @("compName","industry","revenue","employees","status","category","priority","source","type","address","city","state","zip","country")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $cell = $companies.Cells.Item(1, $col + 1)
    $cell.Value = $headers[$col]
    $cell.Interior.Color = 65535
}

$row2 = @("ABC","IT","'3000000","'10000","Active","Client","High","Ad","Billing Address","Mumbai address","Mumbai","Maharashtra","'400028","India")
for ($col = 0; $col -lt $row2.Length; $col++) {
    $companies.Cells.Item(2, $col + 1).Value = $row2[$col]
}

$row3 = @("DEF","FMCG","'140000000","'50000","Inactive","Partner","Medium","Customer","Shipping Address","Thane address","Thane","Maharashtra","'400063","India")
for ($col = 0; $col -lt $row3.Length; $col++) {
    $companies.Cells.Item(3, $col + 1).Value = $row3[$col]
}

# Column widths sized to fit the new data (approximate Excel's "best fit").
# The host pads ColumnWidth by a fixed 0.8333 (5/6) character offset when it
# stores the column, so pre-subtract it to land on the intended widths.
$widths = @(10.44140625,7.44140625,10,9.77734375,7.44140625,8.109375,7.6640625,8.88671875,14.5546875,14.44140625,7.6640625,11.44140625,7,7.21875)
$columnWidthPad = 0.8333333333333333
for ($col = 0; $col -lt $widths.Length; $col++) {
    $companies.Columns.Item($col + 1).ColumnWidth = $widths[$col] - $columnWidthPad
}

$companies.Range("G10").Select() | Out-Null

# Companies becomes the active / visible tab
$companies.Activate()
